$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.604.83"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "3.906.90"
$ws.Range("E3").Value = "  +2.45%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "166.85"
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("D7").Value = "3.906.02"
$ws.Range("E7").Value = "  +2.43%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").Value = "6.45"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "0.0000256"
$ws.Range("E13").Value = "  +3.89%  "
$ws.Range("D14").Value = "37.47"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "4.563.31"
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "3.913.33"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").Value = "68.715.16"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").Value = "17.33"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").Value = "488.93"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "0.727"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("E24").Value = "  +5.48%  "
$ws.Range("D25").Value = "84.69"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").Value = "12.06"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("D28").Value = "10.16"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("D31").Value = "4.061.30"
$ws.Range("E31").Value = "  +2.54%  "
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  -3.70%  "
$ws.Range("D34").Value = "31.89"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("D35").Value = "3.864.19"
$ws.Range("E35").Value = "  +2.75%  "
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").Value = "5.95"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").Value = "3.21"
$ws.Range("E40").Value = "  +5.53%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").Value = "431.11"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D48").Value = "142.66"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").Value = "2.811.57"
$ws.Range("E49").Value = "  -0.79%  "

# Row 50/51 restructure: insert Arweave at 50, shift FLOKI to 51 with updated values, drop VeChain
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").Value = "39.52"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "0.000265"
$ws.Range("E51").Value = "  +16.45%  "
